# "Add files via upload" — append the new menu items that were added to
# the restaurant menu workbook: 19 new rows (33-51) on Sheet1, each with
# an Item name (A), Full price (C) and Image filename (D), using the same
# layout/style as the existing rows. Finish by updating the sheet view to
# match the state Excel leaves behind after the rows were typed in (zoomed
# in a bit, selection parked on the next empty row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-use the currency number format already applied to the existing
# "Full" price column (C32) so the new price cells pick up the same
# style (the "₹ #,##0" format) instead of the plain General format.
$currencyFormat = $ws.Cells.Item(32, 3).NumberFormat

$newItems = @(
    @{ Item = "Thumsup25";                           Price = 25;  Image = "Thusmup Can.jpg" },
    @{ Item = "Chips 20";                             Price = 20;  Image = "Chips 20.jpg" },
    @{ Item = "Chips 10";                             Price = 10;  Image = "Chips 10.jpg" },
    @{ Item = "Chips 05";                             Price = 5;   Image = "Chips 5.jpg" },
    @{ Item = "Ti Tac 5";                             Price = 5;   Image = "Tictac 5.jpg" },
    @{ Item = "Bourbon and Dark Fantasy Biscut 10";   Price = 10;  Image = "Bourbon and Dark Fantasy Biscut 10.jpg" },
    @{ Item = "Cake 15";                              Price = 15;  Image = "britannia cake 15.jpg" },
    @{ Item = "Kitkat 25";                            Price = 25;  Image = "Kitkat 25.jpg" },
    @{ Item = "Cake 1 kg Chocolate";                  Price = 550; Image = "Cake Chocolate.jpg" },
    @{ Item = "Cake 1 kg Pineapple";                  Price = 480; Image = "Cake Pineapple.jpg" },
    @{ Item = "Cake 1kg Butterscotch";                Price = 500; Image = "Cake Butterscotch.jpg" },
    @{ Item = "Cake 600gm Pineapple";                 Price = 350; Image = "Cake Pineapple.jpg" },
    @{ Item = "Cake 600gm Chocolate";                 Price = 450; Image = "Cake Chocolate.jpg" },
    @{ Item = "Cake 600gm Butterscotch";              Price = 400; Image = "Cake Butterscotch.jpg" },
    @{ Item = "Munch 10";                             Price = 10;  Image = "Munch 10.png" },
    @{ Item = "Hide and Seek Black bourbon 10";       Price = 10;  Image = "Hide and Seek Black bourne 10.jpg" },
    @{ Item = "Hide and Seek 10";                     Price = 10;  Image = "Hide and seek 10.jpg" },
    @{ Item = "Hide and Seek 30";                     Price = 30;  Image = "Hide and Seek Biscut 30.jpg" },
    @{ Item = "Thums up 35";                          Price = 35;  Image = "Thums up 35.jpg" }
)

$startRow = 33
for ($i = 0; $i -lt $newItems.Count; $i++) {
    $row = $startRow + $i
    $entry = $newItems[$i]

    $ws.Cells.Item($row, 1).Value = $entry.Item
    $ws.Cells.Item($row, 3).Value = $entry.Price
    $ws.Cells.Item($row, 3).NumberFormat = $currencyFormat
    $ws.Cells.Item($row, 4).Value = $entry.Image
}

# Reflect the post-edit view state: zoomed to 110% with the next empty
# row selected (mirrors what Excel records after the rows were entered).
$ws.Application.ActiveWindow.Zoom = 110
$ws.Range("D52").Select()
